$d = $word.ActiveDocument

# Remove the whole "SMARTREWARDS FAQ'S" Q&A block — from the paragraph
# that begins "How much SmartRewards can I get each month?" through the
# paragraph ending "...every second block 1000 addresses will get paid."
# (The FAQ heading itself, and the trailing empty paragraph after the
# block, are left untouched.)

$startRange = $d.Content.Duplicate
[void]$startRange.Find.Execute("How much SmartRewards can I get each month?", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$startRange.Expand(4)   # wdParagraph -> grow to the full paragraph incl. mark

$endRange = $d.Content.Duplicate
[void]$endRange.Find.Execute("every second block 1000 addresses will get paid.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$endRange.Expand(4)     # wdParagraph -> grow to the full paragraph incl. mark

$deleteRange = $d.Range($startRange.Start, $endRange.End)
[void]$deleteRange.Delete()
